{"js": "// Apply the R-script rename edits described in the commit:\n//   - \"data\"    -> \"houses\"   (the data frame variable)\n//   - \"lin.reg\" -> \"linReg\"   (the lm() model variable)\n//   - \"stdres\"  -> \"mystdres\" (the standardized residuals variable)\n// and merge the \"Example 8: ...\" heading's runs back into a single run\n// (same visible text, just re-typed as one run, matching the target XML).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\nfunction findParagraph(predicate) {\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (predicate(paragraphs.items[i].text)) {\n      return paragraphs.items[i];\n    }\n  }\n  throw new Error(\"No paragraph matched the given predicate\");\n}\n\n// 1) Chapter heading paragraph: \"Example 8:  House Selling Price \u2013 Histogram of\n//    Standardized Residuals\" was re-typed as a single run (text is unchanged).\nconst headingPara = findParagraph((t) => t.indexOf(\"Example \") === 0 && t.indexOf(\"Histogram of Standardized Residuals\") !== -1);\nheadingPara.getRange().insertText(headingPara.text, Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) \"Reading in data\" code chunk: \"data <- read.csv(...)\" -> \"houses <- read.csv(...)\"\nconst readCsvPara = findParagraph((t) => t.indexOf(\"read.csv\") !== -1);\nconst readCsvDataHits = readCsvPara.search(\"data\", { matchCase: true });\nreadCsvDataHits.load(\"text\");\nawait context.sync();\n// The first \"data\" is the variable being assigned; the one inside the URL stays put.\nreadCsvDataHits.items[0].insertText(\"houses\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3) \"Fitting in multiple regression model\" code chunk:\n//    \"lin.reg <- lm(HP.in.thousands ~ House.Size + Bedrooms, data = data)\\rlin.reg\"\n//    -> \"linReg <- lm(HP.in.thousands ~ House.Size + Bedrooms, data = houses)\\rlinReg\"\nconst lmPara = findParagraph((t) => t.indexOf(\"lm(HP.in.thousands\") !== -1);\nconst lmDataHits = lmPara.search(\"data\", { matchCase: true });\nlmDataHits.load(\"text\");\nconst lmLinRegHits = lmPara.search(\"lin.reg\", { matchCase: true });\nlmLinRegHits.load(\"text\");\nawait context.sync();\n// occurrence 0 is the \"data =\" argument name (keep); occurrence 1 is the value (rename).\nlmDataHits.items[1].insertText(\"houses\", Word.InsertLocation.replace);\nfor (let i = 0; i < lmLinRegHits.items.length; i++) {\n  lmLinRegHits.items[i].insertText(\"linReg\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 4) Verbatim model-call echo:\n//    \"## lm(formula = HP.in.thousands ~ House.Size + Bedrooms, data = data)\"\n//    -> \"...data = houses)\"\nconst verbatimPara = findParagraph((t) => t.indexOf(\"## lm(formula = HP.in.thousands\") !== -1);\nconst verbatimDataHits = verbatimPara.search(\"data\", { matchCase: true });\nverbatimDataHits.load(\"text\");\nawait context.sync();\nverbatimDataHits.items[1].insertText(\"houses\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 5) \"To obtain standardized residuals\" code chunk:\n//    \"stdres <- rstandard(lin.reg)\\rhead(stdres)\"\n//    -> \"mystdres <- rstandard(linReg)\\rhead(mystdres)\"\nconst stdresPara = findParagraph((t) => t.indexOf(\"rstandard(\") !== -1);\nconst stdresLinRegHits = stdresPara.search(\"lin.reg\", { matchCase: true });\nstdresLinRegHits.load(\"text\");\nconst stdresHits = stdresPara.search(\"stdres\", { matchCase: true });\nstdresHits.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < stdresLinRegHits.items.length; i++) {\n  stdresLinRegHits.items[i].insertText(\"linReg\", Word.InsertLocation.replace);\n}\nfor (let i = 0; i < stdresHits.items.length; i++) {\n  stdresHits.items[i].insertText(\"mystdres\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 6) Histogram call: \"hist(stdres, breaks = 20, ...\" -> \"hist(mystdres, breaks = 20, ...\"\nconst histPara = findParagraph((t) => t.indexOf(\"hist(\") === 0);\nconst histStdresHits = histPara.search(\"stdres\", { matchCase: true });\nhistStdresHits.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < histStdresHits.items.length; i++) {\n  histStdresHits.items[i].insertText(\"mystdres\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Apply the R-script rename edits described in the commit:\n#   - \"data\"    -> \"houses\"   (the data frame variable)\n#   - \"lin.reg\" -> \"linReg\"   (the lm() model variable)\n#   - \"stdres\"  -> \"mystdres\" (the standardized residuals variable)\n# and re-type the \"Example 8: ...\" heading so its runs collapse back into a\n# single run (same visible text, matching the target XML).\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1) Chapter heading paragraph: \"Example 8:  House Selling Price -\n#    Histogram of Standardized Residuals\" gets re-typed as a single run.\n#    (The shim treats an identical-text write as a no-op, so stage the\n#    change through a throw-away placeholder first to force a real edit.)\n# ---------------------------------------------------------------------\n$paraCount = $d.Paragraphs.Count\nfor ($idx = 1; $idx -le $paraCount; $idx++) {\n    $p = $d.Paragraphs($idx)\n    $t = $p.Range.Text\n    if ($t.StartsWith(\"Example \") -and $t.Contains(\"Histogram of Standardized Residuals\")) {\n        $r = $p.Range\n        $tmpRange = $d.Range($r.Start, $r.End - 1)\n        $tmpRange.Text = \"zzzTEMPzzz\"\n\n        $p2 = $d.Paragraphs($idx)\n        $r2 = $p2.Range\n        $finalRange = $d.Range($r2.Start, $r2.End - 1)\n        $finalRange.Text = \"Example 8:  House Selling Price \" + [char]0x2013 + \" Histogram of Standardized Residuals\"\n        break\n    }\n}\n\n# ---------------------------------------------------------------------\n# 2) Rename \"lin.reg\" -> \"linReg\" everywhere (unambiguous - only refers\n#    to the regression-model variable).\n# ---------------------------------------------------------------------\n$findLinReg = $d.Content.Find\n$findLinReg.ClearFormatting()\n$findLinReg.Text = \"lin.reg\"\n$findLinReg.Replacement.Text = \"linReg\"\n$findLinReg.Execute($null, $true, $true, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# ---------------------------------------------------------------------\n# 3) Rename \"stdres\" -> \"mystdres\" everywhere (unambiguous - only refers\n#    to the standardized-residuals variable).\n# ---------------------------------------------------------------------\n$findStdres = $d.Content.Find\n$findStdres.ClearFormatting()\n$findStdres.Text = \"stdres\"\n$findStdres.Replacement.Text = \"mystdres\"\n$findStdres.Execute($null, $true, $true, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# ---------------------------------------------------------------------\n# 4) Rename \"data\" -> \"houses\", but ONLY the variable-name occurrences:\n#      - \"data <- read.csv(...)\"                       (assignment)\n#      - \"... Bedrooms, data = data)\"                   (the value, 2nd \"data\")\n#      - \"## lm(formula = ... Bedrooms, data = data)\"   (the value, 2nd \"data\")\n#    Leave untouched:\n#      - the \"Reading in data\" heading\n#      - \"data\" inside the CSV URL\n#      - the \"data =\" argument name itself\n# ---------------------------------------------------------------------\n$findData = $d.Content.Find\n$findData.ClearFormatting()\n$findData.Text = \"data\"\n$findData.Execute($null, $true, $true, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null\n$matchNumber = 0\n$occurrencesToRename = @(2, 5, 7)\nwhile ($findData.Found) {\n    $matchNumber = $matchNumber + 1\n    if ($occurrencesToRename -contains $matchNumber) {\n        $findData.Parent.Text = \"houses\"\n    }\n    $findData.Execute($null, $true, $true, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null\n}\n"}
